# This script recreates the edit described by the commit:
# "Added classes Switch and Toggle" - in effect the lab-report document's
# "Контрольные вопросы" (Control questions) section, containing a heading
# and two question/answer pairs, is collapsed into a single paragraph that
# holds a log of simulated circuit-editing actions.

$d = $word.ActiveDocument

# Document originally has 5 paragraphs:
#   1. "Контрольные вопросы:" (Heading1, bold)
#   2. "Вопрос №1"
#   3. "Ответ 1"           <- spacing(after=200,line=310), ind firstLine=855, jc=both
#   4. "Вопрос №2"
#   5. "Ответ 2"           <- spacing(after=200,line=310), ind firstLine=855, jc=both
#
# Target has a single paragraph with spacing(line=310), ind firstLine=855, jc=both,
# and the big replacement text, using the same run formatting (Times New Roman, 14pt,
# black, ru-RU, no bold) as the original "Ответ" runs.

# 1) Remove the heading paragraph ("Контрольные вопросы:") entirely.
$d.Paragraphs(1).Range.Delete()

# 2) Remove the "Вопрос №1" paragraph entirely.
$d.Paragraphs(1).Range.Delete()

# At this point the document is: [1]="Ответ 1" [2]="Вопрос №2" [3]="Ответ 2"

# 3) Merge paragraph 1 ("Ответ 1") into paragraph 2 ("Вопрос №2") by deleting only
#    paragraph 1's end-of-paragraph mark. The merged paragraph keeps paragraph 2's
#    (clean, no jc/after) paragraph properties, while its first run keeps "Ответ 1"'s
#    (non-bold) run formatting.
$p1 = $d.Paragraphs(1)
$d.Range($p1.Range.End - 1, $p1.Range.End).Delete()

# Now: [1]="Ответ 1Вопрос №2" (clean pPr) [2]="Ответ 2"

# 4) Remove the trailing "Ответ 2" paragraph entirely.
$d.Paragraphs(2).Range.Delete()

# Now a single paragraph remains: "Ответ 1Вопрос №2"

# 5) Replace the whole paragraph text with the final combined text via Find/Replace
#    (robust across the paragraph's multiple runs), keeping the first run's formatting.
$p = $d.Paragraphs(1)
$rng = $p.Range
$rng.Find.Execute($rng.Text.TrimEnd([char]13), $true, $false, $false, $false, $false, $true, 1, $false, "На схему был добавлен амперметр. Показания на амперметре были изменены до 20 A. На схему был добавлен вольтметр. Показания на амперметре были изменены до 20 A. Показания на вольтметре были изменены до 24 В. На схему был добавлен мультиметр. На схему был добавлен мультиметр. На схему был добавлен мультиметр. На схему был добавлен резистор. Сопротивление резистора было увеличено до 1 Ом. Сопротивление резистора было увеличено до 2 Ом. Сопротивление резистора было увеличено до 3 Ом. Сопротивление резистора было уменьшено до 2 Ом. На схему был добавлен проводник. Длина проводника изменена до 1 cм, также диаметр изменён до 2 см, и удельное электрическое сопротивление материала проводника изменено до 3. Сопротивление проводника было изменено до 0,95 Ом. На схему был добавлен реостат. Сопротивление реостата было изменено до 38 Ом. На схему был добавлен источник напряжения. На схему был добавлен источник напряжения. На схему был добавлен конденсатор. Площадь пластин плоского конденсатора изменена до 12 cм^2, также значение относительной диэлектрической проницаемости изменено до 34, и расстояние между пластинами изменено до 4 мм. На схему был добавлен амперметр. Внутренний и внешний радиусы цилиндрического конденсатора изменены до 3 см и 4 см соответственно, также высота конденсатора изменена до 6 см, и значение относительной диэлектрической проницаемости изменено до 56. Внутренний и внешний радиусы цилиндрического конденсатора изменены до 3 см и 4 см соответственно, также высота конденсатора изменена до 6 см, и значение относительной диэлектрической проницаемости изменено до 56. Показания на амперметре были изменены до 20 A. Показания на вольтметре были изменены до 24 В. Сопротивление реостата было изменено до 18 Ом. ", 2)

# 6) Justify the paragraph (w:jc w:val="both").
$d.Paragraphs(1).Format.Alignment = 3
